# Generate Report for Archive
# - flips the stale "Ready for handoff" status to "In Translation"
#   everywhere it appears (Overview summary columns + per-locale Status
#   columns), and re-fits the affected "Status" columns to their new,
#   shorter content.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- 1. Replace the status text wherever it occurs on every sheet ---
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $firstRow = $used.Row
    $firstCol = $used.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($firstRow + $r, $firstCol + $c)
            if ($oldStatus -eq $cell.Value2) {
                $cell.Value = $newStatus
            }
        }
    }
}

# --- 2. Re-fit the Status columns now that the text is shorter ---
# Overview sheet keeps one status column per locale (zh-cn -> E, de-de -> F)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# Each locale sheet (zh-cn, de-de) has its own "Status" column (column C)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
